$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "60.796.62"
Set-TextValue "E2" "  -1.70%  "
Set-TextValue "D3" "3.387.51"
Set-TextValue "E3" "  -2.10%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "571.51"
Set-TextValue "E5" "  -1.58%  "
Set-TextValue "D6" "141.03"
Set-TextValue "E6" "  -4.84%  "
Set-TextValue "E7" "  +0.10%  "
Set-TextValue "D8" "3.388.36"
Set-TextValue "E8" "  -2.11%  "
Set-TextValue "E9" "  -0.46%  "
Set-TextValue "D10" "7.49"
Set-TextValue "E10" "  -2.93%  "
Set-TextValue "D11" "0.123"
Set-TextValue "E11" "  -1.29%  "
Set-TextValue "D12" "0.391"
Set-TextValue "E12" "  +0.18%  "
Set-TextValue "D13" "3.968.17"
Set-TextValue "E13" "  -2.09%  "
Set-TextValue "E14" "  +0.39%  "
Set-TextValue "E15" "  +0.40%  "
Set-TextValue "D16" "0.0000170"
Set-TextValue "E16" "  -2.97%  "
Set-TextValue "D17" "3.390.27"
Set-TextValue "E17" "  -2.35%  "
Set-TextValue "D18" "60.941.97"
Set-TextValue "E18" "  -1.50%  "
Set-TextValue "D19" "6.26"
Set-TextValue "E19" "  -1.38%  "
Set-TextValue "D20" "14.11"
Set-TextValue "E20" "  -2.35%  "
Set-TextValue "D21" "9.03"
Set-TextValue "E21" "  -4.53%  "
Set-TextValue "D22" "387.29"
Set-TextValue "E22" "  +0.35%  "
Set-TextValue "D23" "0.558"
Set-TextValue "E23" "  -2.11%  "
Set-TextValue "E24" "  +0.73%  "
Set-TextValue "D25" "1.00"
Set-TextValue "E25" "  -0.17%  "
Set-TextValue "D26" "0.0000118"
Set-TextValue "E26" "  -3.84%  "
Set-TextValue "D27" "3.530.32"
Set-TextValue "E27" "  -1.78%  "
Set-TextValue "E28" "  -0.82%  "
Set-TextValue "D29" "0.999"
Set-TextValue "E29" "  +0.08%  "
Set-TextValue "D30" "7.36"
Set-TextValue "E30" "  -5.84%  "
Set-TextValue "D31" "8.12"
Set-TextValue "E31" "  -1.80%  "
Set-TextValue "B32" "PancakeSwap"
Set-TextValue "C32" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D32" "2.15"
Set-TextValue "E32" "  -1.24%  "
Set-TextValue "B33" "Fetch.AI"
Set-TextValue "C33" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D33" "1.44"
Set-TextValue "E33" "  -7.50%  "
Set-TextValue "E34" "  -0.01%  "
Set-TextValue "D35" "23.77"
Set-TextValue "E35" "  -1.03%  "
Set-TextValue "B36" "RenzoRestakedETH"
Set-TextValue "C36" "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue "D36" "3.416.95"
Set-TextValue "E36" "  -1.89%  "
Set-TextValue "B37" "Aptos"
Set-TextValue "C37" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D37" "6.91"
Set-TextValue "E37" "  -2.06%  "
Set-TextValue "D38" "167.03"
Set-TextValue "E38" "  -0.09%  "
Set-TextValue "D39" "5.04"
Set-TextValue "E39" "  -3.32%  "
Set-TextValue "D40" "1.51"
Set-TextValue "E40" "  -2.86%  "
Set-TextValue "D41" "0.0778"
Set-TextValue "E41" "  -1.56%  "
Set-TextValue "D42" "26.68"
Set-TextValue "E42" "  +2.72%  "
Set-TextValue "E43" "  -1.36%  "
Set-TextValue "E44" "  +0.12%  "
Set-TextValue "D45" "4.44"
Set-TextValue "E45" "  -1.18%  "
Set-TextValue "D46" "41.81"
Set-TextValue "E46" "  -1.28%  "
Set-TextValue "E47" "  -2.69%  "
Set-TextValue "D48" "2.531.39"
Set-TextValue "E48" "  -2.88%  "
Set-TextValue "D49" "1.12"
Set-TextValue "E49" "  -4.28%  "
Set-TextValue "D50" "6.84"
Set-TextValue "E50" "  -1.84%  "
Set-TextValue "D51" "22.82"
Set-TextValue "E51" "  -2.63%  "
